$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 311, shifting existing rows (311-358) down to (312-359).
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row 311 with the same "template" values as the row
# that is now directly below it (old row 311, now row 312), then overwrite the
# columns that actually carry new data for this entry.
for ($col = 1; $col -le 18; $col++) {
    $srcCell = $ws.Cells.Item(312, $col)
    $dstCell = $ws.Cells.Item(311, $col)
    $dstCell.Value2 = $srcCell.Value2
}

# Now set the specific new values for this new record.
$ws.Cells.Item(311, 4).Value2 = 45180   # D311 Fecha
$ws.Cells.Item(311, 10).Value2 = 220    # J311 Volumen
$ws.Cells.Item(311, 13).Value2 = 11091  # M311 Precio promedio ponderado
$ws.Cells.Item(311, 16).Value2 = 222    # P311 Precio $/Kg
